$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A26").Value = "BIOMI2500"
$ws.Range("A27").Value = "BIOMS3310"

$ws.Range("A28").Select()
